# Automatic update of files.
# Inserts two new observation records (rows 4 and 5) into the "Artfynd" sheet,
# pushing the previous rows 4-10 down to rows 6-12, and updates the
# Taxonsorteringsordning values on rows 2 and 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2 and 3 (only column B changes) ---
$ws.Range("B2").Value2 = 90804
$ws.Range("B3").Value2 = 88623

# --- Insert two blank rows at position 4, pushing old rows 4-10 to 6-12 ---
$ws.Rows("4:5").Insert()

# --- Populate new row 4 ---
$ws.Range("A4").Value2 = 112472885
$ws.Range("B4").Value2 = 90804
$ws.Range("C4").Value2 = "Ovaliderad"
$ws.Range("D4").Value2 = "VU"
$ws.Range("E4").Value2 = 4365
$ws.Range("F4").Value2 = "Smalfotad taggsvamp"
$ws.Range("G4").Value2 = "Hydnellum gracilipes"
$ws.Range("H4").Value2 = "(P.Karst) P.Karst"
$ws.Range("P4").Value2 = "Prästtjärnen (Prästtjärnen), Dlr"
$ws.Range("Q4").Value2 = 517956
$ws.Range("R4").Value2 = 6790407
$ws.Range("S4").Value2 = 25
$ws.Range("T4").Value2 = "Dalarna"
$ws.Range("U4").Value2 = "Rättvik"
$ws.Range("V4").Value2 = "Dalarna"
$ws.Range("W4").Value2 = "Ore"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value2 = "2023-10-02"
$ws.Range("Z4").Value2 = "13:12"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value2 = "2023-10-02"
$ws.Range("AB4").Value2 = "13:12"
$ws.Range("AD4").Value2 = $false
$ws.Range("AE4").Value2 = $false
$ws.Range("AG4").Value2 = $false
$ws.Range("AW4").Value2 = "Andreas Öster"
$ws.Range("AX4").Value2 = "Andreas Öster"

# --- Populate new row 5 ---
$ws.Range("A5").Value2 = 112473083
$ws.Range("B5").Value2 = 89806
$ws.Range("C5").Value2 = "Ovaliderad"
$ws.Range("D5").Value2 = "EN"
$ws.Range("E5").Value2 = 71
$ws.Range("F5").Value2 = "Urskogsporing"
$ws.Range("G5").Value2 = "Neoantrodia infirma"
$ws.Range("H5").Value2 = "(Renvall & Niemelä) Audet"
$ws.Range("P5").Value2 = "Prästtjärnen (Prästtjärnen), Dlr"
$ws.Range("Q5").Value2 = 518039
$ws.Range("R5").Value2 = 6790377
$ws.Range("S5").Value2 = 25
$ws.Range("T5").Value2 = "Dalarna"
$ws.Range("U5").Value2 = "Rättvik"
$ws.Range("V5").Value2 = "Dalarna"
$ws.Range("W5").Value2 = "Ore"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value2 = "2023-10-02"
$ws.Range("Z5").Value2 = "13:17"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value2 = "2023-10-02"
$ws.Range("AB5").Value2 = "13:17"
$ws.Range("AD5").Value2 = $false
$ws.Range("AE5").Value2 = $false
$ws.Range("AG5").Value2 = $false
$ws.Range("AW5").Value2 = "Andreas Öster"
$ws.Range("AX5").Value2 = "Andreas Öster"

Write-Host "Edit complete"
